$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 (bold, centered, bordered) onto the two
# new header cells so they pick up the same cellXf (s="1") instead of a
# freshly minted style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data columns I (I0) and J (IF) for rows 2-12
$data = @(
    @(7, 7),
    @(6, 6),
    @(6, 6),
    @(6, 8),
    @(4, 4),
    @(5, 6),
    @(5, 6),
    @(1, 2),
    @(8, 8),
    @(1, 3),
    @(1, 2)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
